$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B3").Value = 27698
$ws.Range("C3").Value = 4547
$ws.Range("D3").Value = 4625
$ws.Range("B4").Value = 16367
$ws.Range("C4").Value = 1737
$ws.Range("D4").Value = 1982
$ws.Range("B5").Value = 45788
$ws.Range("C5").Value = 4495
$ws.Range("D5").Value = 4343
$ws.Range("B6").Value = 812
$ws.Range("C6").Value = 441
$ws.Range("D6").Value = 101
$ws.Range("B7").Value = 30886
$ws.Range("C7").Value = 5285
$ws.Range("D7").Value = 4046
$ws.Range("B8").Value = 3818
$ws.Range("C8").Value = 652
$ws.Range("D8").Value = 910
$ws.Range("B9").Value = 3506
$ws.Range("C9").Value = 761
$ws.Range("D9").Value = 502
$ws.Range("B10").Value = 1680
$ws.Range("C10").Value = 278
$ws.Range("D10").Value = 158
$ws.Range("B11").Value = 395
$ws.Range("C11").Value = 209
$ws.Range("D11").Value = 2
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 0
$ws.Range("B13").Value = 660
$ws.Range("C13").Value = 134
$ws.Range("D13").Value = 189
$ws.Range("B14").Value = 1985
$ws.Range("C14").Value = 753
$ws.Range("D14").Value = 810
$ws.Range("B15").Value = 3219
$ws.Range("C15").Value = 1151
$ws.Range("D15").Value = 611
$ws.Range("B16").Value = 2335
$ws.Range("C16").Value = 1099
$ws.Range("D16").Value = 268
$ws.Range("B17").Value = 1803
$ws.Range("C17").Value = 516
$ws.Range("D17").Value = 237
$ws.Range("B18").Value = 11548
$ws.Range("C18").Value = 1625
$ws.Range("D18").Value = 2154
$ws.Range("B19").Value = 836
$ws.Range("C19").Value = 421
$ws.Range("D19").Value = 234
$ws.Range("B20").Value = 11518
$ws.Range("C20").Value = 1731
$ws.Range("D20").Value = 1931
$ws.Range("B21").Value = 159
$ws.Range("C21").Value = 293
$ws.Range("D21").Value = 11
$ws.Range("B22").Value = 10604
$ws.Range("C22").Value = 1382
$ws.Range("D22").Value = 1761
$ws.Range("B23").Value = 665
$ws.Range("C23").Value = 481
$ws.Range("D23").Value = 114
$ws.Range("B24").Value = 12117
$ws.Range("C24").Value = 1296
$ws.Range("D24").Value = 2434
$ws.Range("B25").Value = 51822
$ws.Range("C25").Value = 4642
$ws.Range("D25").Value = 6262
$ws.Range("B26").Value = 3897
$ws.Range("C26").Value = 1243
$ws.Range("D26").Value = 560
$ws.Range("B27").Value = 0
$ws.Range("C27").Value = 0
$ws.Range("D27").Value = 0
$ws.Range("B28").Value = 3198
$ws.Range("C28").Value = 802
$ws.Range("D28").Value = 775
$ws.Range("B29").Value = 1216
$ws.Range("C29").Value = 373
$ws.Range("D29").Value = 281
$ws.Range("B30").Value = 9114
$ws.Range("C30").Value = 1634
$ws.Range("D30").Value = 1912
$ws.Range("B31").Value = 341
$ws.Range("C31").Value = 79
$ws.Range("D31").Value = 197
$ws.Range("B32").Value = 1931
$ws.Range("C32").Value = 1166
$ws.Range("D32").Value = 199
$ws.Range("B33").Value = 8849
$ws.Range("C33").Value = 2156
$ws.Range("D33").Value = 1945
$ws.Range("B34").Value = 6404
$ws.Range("C34").Value = 2241
$ws.Range("D34").Value = 1213
$ws.Range("B35").Value = 3139
$ws.Range("C35").Value = 366
$ws.Range("D35").Value = 628
$ws.Range("B36").Value = 36459
$ws.Range("C36").Value = 4092
$ws.Range("D36").Value = 3469
$ws.Range("B37").Value = 5100
$ws.Range("C37").Value = 1825
$ws.Range("D37").Value = 733
$ws.Range("B38").Value = 17371
$ws.Range("C38").Value = 1178
$ws.Range("D38").Value = 1661
$ws.Range("B39").Value = 487
$ws.Range("C39").Value = 727
$ws.Range("D39").Value = 119
$ws.Range("B40").Value = 798
$ws.Range("C40").Value = 349
$ws.Range("D40").Value = 365
$ws.Range("B41").Value = 1938
$ws.Range("C41").Value = 447
$ws.Range("D41").Value = 68
$ws.Range("B42").Value = 7403
$ws.Range("C42").Value = 347
$ws.Range("D42").Value = 268
$ws.Range("B43").Value = 177
$ws.Range("C43").Value = 148
$ws.Range("D43").Value = 13
$ws.Range("B44").Value = 514
$ws.Range("C44").Value = 111
$ws.Range("D44").Value = 49
$ws.Range("B45").Value = 0
$ws.Range("C45").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("B46").Value = 2002
$ws.Range("C46").Value = 738
$ws.Range("D46").Value = 270
$ws.Range("B47").Value = 8814
$ws.Range("C47").Value = 2298
$ws.Range("D47").Value = 1682
$ws.Range("B48").Value = 21886
$ws.Range("C48").Value = 2274
$ws.Range("D48").Value = 3462
$ws.Range("B49").Value = 9666
$ws.Range("C49").Value = 2200
$ws.Range("D49").Value = 1043
$ws.Range("B50").Value = 7538
$ws.Range("C50").Value = 994
$ws.Range("D50").Value = 1652
$ws.Range("B51").Value = 21468
$ws.Range("C51").Value = 1920
$ws.Range("D51").Value = 3577
$ws.Range("B52").Value = 3015
$ws.Range("C52").Value = 705
$ws.Range("D52").Value = 769
$ws.Range("B53").Value = 7313
$ws.Range("C53").Value = 1459
$ws.Range("D53").Value = 1225
$ws.Range("B54").Value = 1303
$ws.Range("C54").Value = 889
$ws.Range("D54").Value = 532
$ws.Range("B55").Value = 1469
$ws.Range("C55").Value = 936
$ws.Range("D55").Value = 73
$ws.Range("B56").Value = 3277
$ws.Range("C56").Value = 808
$ws.Range("D56").Value = 1204
$ws.Range("B57").Value = 8199
$ws.Range("C57").Value = 3137
$ws.Range("D57").Value = 1878
$ws.Range("B58").Value = 10032
$ws.Range("C58").Value = 1053
$ws.Range("D58").Value = 398
$ws.Range("B59").Value = 436308
$ws.Range("C59").Value = 68434
$ws.Range("D59").Value = 63629
